$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") '330.67'
Set-TextValue $ws.Range("E2") '0.05%'

# Row 3
Set-TextValue $ws.Range("D3") '45.39'
Set-TextValue $ws.Range("E3") '2.22%'

# Row 4
Set-TextValue $ws.Range("D4") '5.578'
Set-TextValue $ws.Range("E4") '1.59%'

# Row 5
Set-TextValue $ws.Range("D5") '0.08346'
Set-TextValue $ws.Range("E5") '3.94%'

# Row 6
Set-TextValue $ws.Range("D6") '2.096'
Set-TextValue $ws.Range("E6") '0.05%'

# Row 7
Set-TextValue $ws.Range("D7") '0.9802'
Set-TextValue $ws.Range("E7") '2.90%'

# Row 8
Set-TextValue $ws.Range("D8") '2.543'
Set-TextValue $ws.Range("E8") '-0.65%'

# Row 9
Set-TextValue $ws.Range("D9") '0.1204'
Set-TextValue $ws.Range("E9") '5.16%'

# Row 10
Set-TextValue $ws.Range("D10") '0.1915'
Set-TextValue $ws.Range("E10") '1.40%'

# Row 11
Set-TextValue $ws.Range("D11") '10.30'
Set-TextValue $ws.Range("E11") '-3.74%'

# Row 12
Set-TextValue $ws.Range("D12") '0.09863'
Set-TextValue $ws.Range("E12") '0.17%'

# Row 13
Set-TextValue $ws.Range("D13") '0.04642'
Set-TextValue $ws.Range("E13") '-3.81%'

# Row 14
Set-TextValue $ws.Range("E14") '-0.75%'

# Row 15
Set-TextValue $ws.Range("D15") '0.001291'
Set-TextValue $ws.Range("E15") '1.18%'

# Row 16
Set-TextValue $ws.Range("D16") '0.005895'
Set-TextValue $ws.Range("E16") '-0.30%'

# Row 17
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue $ws.Range("D17") '3.375'
Set-TextValue $ws.Range("E17") '0.25%'

# Row 18
$ws.Range("B18").Value = 'GateToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue $ws.Range("D18") '4.446'
Set-TextValue $ws.Range("E18") '1.00%'

# Row 19
$ws.Range("B19").Value = 'BitpandaEcosystemToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
Set-TextValue $ws.Range("D19") '0.3341'
Set-TextValue $ws.Range("E19") '-3.45%'

# Row 20
$ws.Range("B20").Value = 'ProBitToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
Set-TextValue $ws.Range("D20") '0.1391'
Set-TextValue $ws.Range("E20") '-0.58%'

# Row 21
$ws.Range("B21").Value = 'ZBToken'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
Set-TextValue $ws.Range("D21") '0.2784'
Set-TextValue $ws.Range("E21") '11.23%'

# Row 22
$ws.Range("B22").Value = 'CoinExToken'
$ws.Range("C22").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
Set-TextValue $ws.Range("D22") '0.04176'
Set-TextValue $ws.Range("E22") '2.42%'

# Row 23
$ws.Range("B23").Value = 'BitKan'
$ws.Range("C23").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
Set-TextValue $ws.Range("D23") '0.001293'
Set-TextValue $ws.Range("E23") '1.53%'

# Row 24
$ws.Range("B24").Value = 'HotbitToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
Set-TextValue $ws.Range("D24") '0.004579'
Set-TextValue $ws.Range("E24") '5.01%'

# Row 25
Set-TextValue $ws.Range("D25") '0.0001303'
Set-TextValue $ws.Range("E25") '8.60%'

# Row 26
Set-TextValue $ws.Range("D26") '0.0003745'
Set-TextValue $ws.Range("E26") '0.06%'

# Row 38
Set-TextValue $ws.Range("D38") '0.02697'
Set-TextValue $ws.Range("E38") '3.93%'

# Row 39
Set-TextValue $ws.Range("D39") '0.05754'
Set-TextValue $ws.Range("E39") '-1.30%'

# Row 40
Set-TextValue $ws.Range("D40") '0.007909'
Set-TextValue $ws.Range("E40") '4.66%'

# Row 41
Set-TextValue $ws.Range("D41") '0.1432'
Set-TextValue $ws.Range("E41") '1.98%'

# Row 42
Set-TextValue $ws.Range("D42") '0.007516'
Set-TextValue $ws.Range("E42") '5.17%'

# Row 44
Set-TextValue $ws.Range("D44") '0.008500'
Set-TextValue $ws.Range("E44") '-3.50%'

# Row 45
Set-TextValue $ws.Range("D45") '0.3368'

# Row 46
Set-TextValue $ws.Range("D46") '0.00007120'
Set-TextValue $ws.Range("E46") '1.99%'

# Row 47
Set-TextValue $ws.Range("E47") '0.19%'

# Row 48
Set-TextValue $ws.Range("E48") '0.33%'

# Row 49
$ws.Range("B49").Value = 'CoinbaseStockToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
Set-TextValue $ws.Range("D49") '0.003535'
Set-TextValue $ws.Range("E49") '0.17%'

# Row 50
$ws.Range("B50").Value = 'BOLO'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
Set-TextValue $ws.Range("D50") '0.003530'
Set-TextValue $ws.Range("E50") '-0.54%'

# Row 51
Set-TextValue $ws.Range("D51") '0.00002103'
Set-TextValue $ws.Range("E51") '0.19%'
